$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the match data (columns F:V) between rows 88 and 89 ---
# (columns A:E -- index/pais/torneio/temporada/data_partida -- are left untouched)
for ($col = 6; $col -le 22; $col++) {
    $v88 = $ws.Cells.Item(88, $col).Value2
    $v89 = $ws.Cells.Item(89, $col).Value2
    $ws.Cells.Item(88, $col).Value = $v89
    $ws.Cells.Item(89, $col).Value = $v88
}

# --- Append a new row 139 for the Girona x Ath Bilbao match ---
# Clone formatting from the last existing row (138), then fill the values.
$ws.Range("A138:V138").Copy()
$ws.Range("A139:V139").PasteSpecial(-4122)

$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = "spain"
$ws.Cells.Item(139, 3).Value = "laliga"
$ws.Cells.Item(139, 4).Value = "2023-2024"
$ws.Cells.Item(139, 5).Value = 45257.875
$ws.Cells.Item(139, 6).Value = "Girona"
$ws.Cells.Item(139, 7).Value = 1
$ws.Cells.Item(139, 8).Value = "Ath Bilbao"
$ws.Cells.Item(139, 9).Value = 1
$ws.Cells.Item(139, 10).Value = 2.61
$ws.Cells.Item(139, 11).Value = "05/11/2023 11:03"
$ws.Cells.Item(139, 12).Value = 2.27
$ws.Cells.Item(139, 13).Value = "27/11/2023 20:56"
$ws.Cells.Item(139, 14).Value = 3.31
$ws.Cells.Item(139, 15).Value = "05/11/2023 11:03"
$ws.Cells.Item(139, 16).Value = 3.66
$ws.Cells.Item(139, 17).Value = "27/11/2023 20:59"
$ws.Cells.Item(139, 18).Value = 2.82
$ws.Cells.Item(139, 19).Value = "05/11/2023 11:03"
$ws.Cells.Item(139, 20).Value = 3.17
$ws.Cells.Item(139, 21).Value = "27/11/2023 20:59"
$ws.Cells.Item(139, 22).Value = "https://www.betexplorer.com/football/spain/laliga/girona-ath-bilbao/OnDpwYRA/"
